$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (col D) and Volume(1h) (col E) updates. D values are prefixed with an
# apostrophe to force text interpretation (several look like numbers/dates,
# e.g. "206.47" or "0.0520", and Excel would otherwise coerce them), then the
# style is normalized back so no stray numeric/quote-prefix formatting sticks
# around on the cell.
$updates = @(
    @{Row=2;  D="26.924.31";  E="  -0.40%  "},
    @{Row=3;  D="1.552.30";   E="  -0.47%  "},
    @{Row=4;  D=$null;        E="  -0.43%  "},
    @{Row=5;  D="206.47";     E="  -0.48%  "},
    @{Row=6;  D="0.489";      E="  +0.49%  "},
    @{Row=7;  D=$null;        E="  -0.42%  "},
    @{Row=8;  D=$null;        E="  +0.99%  "},
    @{Row=9;  D=$null;        E="  -0.29%  "},
    @{Row=10; D=$null;        E="  +0.46%  "},
    @{Row=11; D=$null;        E="  -0.78%  "},
    @{Row=12; D="1.772.86";   E="  -0.43%  "},
    @{Row=13; D="1.554.08";   E="  -0.30%  "},
    @{Row=14; D=$null;        E="  +0.45%  "},
    @{Row=15; D=$null;        E="  +0.22%  "},
    @{Row=16; D="26.914.75";  E="  -0.40%  "},
    @{Row=17; D="61.63";      E="  -0.69%  "},
    @{Row=18; D="0.0₃0713";   E="  +3.41%  "},
    @{Row=19; D="216.93";     E="  +0.33%  "},
    @{Row=20; D="7.30";       E="  +0.19%  "},
    @{Row=22; D="4.09";       E="  +1.09%  "},
    @{Row=23; D=$null;        E="  -0.53%  "},
    @{Row=24; D=$null;        E="  -1.17%  "},
    @{Row=25; D="153.82";     E="  +0.61%  "},
    @{Row=26; D="6.64";       E="  -0.53%  "},
    @{Row=27; D="14.97";      E="  +0.10%  "},
    @{Row=28; D=$null;        E="  +0.35%  "},
    @{Row=29; D=$null;        E="  -0.18%  "},
    @{Row=30; D=$null;        E="  +0.96%  "},
    @{Row=31; D=$null;        E="  -1.20%  "},
    @{Row=32; D=$null;        E="  -0.48%  "},
    @{Row=33; D=$null;        E="  +3.52%  "},
    @{Row=34; D="1.411.03";   E="  +0.45%  "},
    @{Row=35; D=$null;        E="  +1.85%  "},
    @{Row=36; D="0.964";      E="  -0.15%  "},
    @{Row=37; D=$null;        E="  +0.32%  "},
    @{Row=38; D=$null;        E="  +0.08%  "},
    @{Row=39; D=$null;        E="  +0.03%  "},
    @{Row=40; D=$null;        E="  -0.53%  "},
    @{Row=41; D=$null;        E="  -0.40%  "},
    @{Row=42; D=$null;        E="  +3.14%  "},
    @{Row=43; D="2.31";       E="  +1.57%  "},
    @{Row=44; D=$null;        E="  +0.44%  "},
    @{Row=45; D=$null;        E="  +0.58%  "},
    @{Row=46; D=$null;        E="  -1.41%  "},
    @{Row=47; D="1.687.01";   E="  -0.43%  "},
    @{Row=48; D="87.34";      E="  +1.22%  "},
    @{Row=49; D="0.0520";     E="  +1.60%  "},
    @{Row=50; D=$null;        E="  +2.92%  "},
    @{Row=51; D="0.0959";     E="  -0.30%  "}
)

foreach ($u in $updates) {
    if ($u.D -ne $null) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
